# dsl/src/test/data/TabularActivityBuilderAndSplitDef.xlsx
# commit: "build splits other than AdnSplit #422"
#
# Rework the "AndSplitWithOrSplit" sheet's OrSplit-flavoured activity names
# into the generic "Split" naming (TestItem_Split:0 / EndSplit / SplitLeft /
# SplitRight), and move the active selection around a couple of sheets to
# match the author's last-saved cursor position.

$wb = $excel.ActiveWorkbook

# --- AndSplitWithLoop: just moves the saved selection, no content change ---
$wsLoop = $wb.Worksheets.Item("AndSplitWithLoop")
$wsLoop.Activate()
$wsLoop.Range("A9").Select()

# --- AndSplitWithOrSplit: rename the Split-related activity references ---
$wsSplit = $wb.Worksheets.Item("AndSplitWithOrSplit")

# Order matters: new shared-string entries are appended in first-use order,
# and the target file expects them as 37=TestItem_Split:0, 38=EndSplit,
# 39=SplitLeft, 40=SplitRight.
$wsSplit.Range("B6").Value = "TestItem_Split:0"
$wsSplit.Range("A11").Value = "EndSplit"
$wsSplit.Range("C6").Value = "SplitLeft"
$wsSplit.Range("C9").Value = "SplitRight"
$wsSplit.Range("B9").Value = "TestItem_Split:0"
$wsSplit.Range("B13").Value = "TestItem_Split:0"

# This sheet becomes the active tab, with C10 as the last selected cell
# (this also drops tabSelected from whichever sheet had it before).
$wsSplit.Activate()
$wsSplit.Range("C10").Select()
